# The deck ships two theme parts:
#   ppt/theme/theme1.xml -> bound to the (one and only) slide master, name="Integral" / clrScheme "Red Violet"
#   ppt/theme/theme2.xml -> bound only to the notes master,           name="Office Theme" / clrScheme "Office"
#
# The commit swaps the two themes' contents wholesale: theme1.xml becomes the
# "Office Theme" palette and theme2.xml becomes the "Integral"/"Red Violet"
# palette. Font scheme + format scheme are already identical between the two
# theme parts, so the only observable content delta is the 12-slot colour
# scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink).
#
# The slide master's live theme (theme1.xml) is reachable and editable through
# the ThemeColorScheme object (1-based, same order as the OOXML clrScheme
# children: dk1, lt1, dk2, lt2, accent1..accent6, hlink, folHlink). PowerPoint's
# ColorFormat.RGB takes/returns a COLORREF (0x00BBGGRR), i.e. R + G*256 + B*65536.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Target palette for theme1.xml after the edit: the stock Office theme colours.
$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = $officeColors[$i - 1]
}
